$wb = $excel.ActiveWorkbook

# The existing "总计" sheet is currently the 4th sheet. We insert a new
# "2022-Q1" sheet right before it (so it keeps its place right after
# "2021-Q3"), and then add a new "2022-Q1" summary row at the top of "总计".
$totalSheet = $wb.Worksheets.Item("总计")

# Duplicate "总计" to use as a formatting template (styles/borders carried
# over) for the new per-fund holdings sheet, then rename + fill it with
# data. NOTE: this shifts "总计" from index 4 to index 5, so after this we
# must always re-fetch sheets by name (index-based sheet variables do not
# track the sheet across insertions in this environment).
$totalSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# Force text-typed columns to avoid Excel auto-coercing numeric-looking
# strings (e.g. leading-zero fund codes) into numbers, then restore the
# default "Normal" style afterwards so no stray number-format style is
# left behind on these cells.
$textRange = $newSheet.Range("B1:G5")
$textRange.NumberFormat = "@"

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Row 2
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "010761"
$newSheet.Cells.Item(2, 3).Value = "华商甄选回报混合"
$newSheet.Cells.Item(2, 4).Value = "20.63"
$newSheet.Cells.Item(2, 5).Value = "93.93"
$newSheet.Cells.Item(2, 6).Value = "5.11"
$newSheet.Cells.Item(2, 7).Value = "1.0542"
$newSheet.Cells.Item(2, 8).Value = 3

# Row 3
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "241001"
$newSheet.Cells.Item(3, 3).Value = "华宝海外中国混合(QDII)"
$newSheet.Cells.Item(3, 4).Value = "0.83"
$newSheet.Cells.Item(3, 5).Value = "86.89"
$newSheet.Cells.Item(3, 6).Value = "4.26"
$newSheet.Cells.Item(3, 7).Value = "0.0354"
$newSheet.Cells.Item(3, 8).Value = 8

# Row 4
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "161124"
$newSheet.Cells.Item(4, 3).Value = "易方达香港恒生综合小型股指数（QDII-LOF）A"
$newSheet.Cells.Item(4, 4).Value = "0.28"
$newSheet.Cells.Item(4, 5).Value = "92.62"
$newSheet.Cells.Item(4, 6).Value = "3.21"
$newSheet.Cells.Item(4, 7).Value = "0.0090"
$newSheet.Cells.Item(4, 8).Value = 1

# Row 5 (brand-new row added at the bottom of the fund table)
$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "006263"
$newSheet.Cells.Item(5, 3).Value = "易方达香港恒生综合小型股指数（QDII-LOF）C"
$newSheet.Cells.Item(5, 4).Value = "0.06"
$newSheet.Cells.Item(5, 5).Value = "92.62"
$newSheet.Cells.Item(5, 6).Value = "3.21"
$newSheet.Cells.Item(5, 7).Value = "0.0019"
$newSheet.Cells.Item(5, 8).Value = 1

$textRange.Style = "Normal"
# Re-apply the header/index-column styles that "Normal" wiped out, copying
# them (via copy/paste-special of formats) from the (freshly re-fetched)
# "总计" sheet that was used as template.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122) # xlPasteFormats

# --- Now update the "总计" sheet: add a new row for 2022-Q1 at the top of
# the data (row 2), pushing the existing rows down to rows 3-5. We give the
# new row 2 the same formatting as the (old) row 2 by copying its format
# first, then rewrite every data row's values directly (rows 2-5) with
# their final, shifted-down contents, which avoids relying on any
# particular Insert/shift semantics.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A5:D5").PasteSpecial(-4122) # xlPasteFormats - row 5 gets row 2's look

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q1"
$totalSheet.Cells.Item(5, 3).Value = 8
$totalSheet.Cells.Item(5, 4).Value = 2.51

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(4, 3).Value = 7
$totalSheet.Cells.Item(4, 4).Value = 5.61

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(3, 3).Value = 7
$totalSheet.Cells.Item(3, 4).Value = 6.07

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 1.1

$wb.Worksheets.Item(1).Activate()
